# Add a new "ExecuteStatus" column (col G) to the FreeCRM_AddContact sheet,
# shifting the existing "TestResult" column to H, and populate it with Y/N
# per-row, plus set keyword "close" on row 28 (Keyword column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FreeCRM_AddContact")

# --- Shift TestResult header + column out of the way: insert a new column at G ---
$ws.Range("G1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("G2").Value = "ExecuteStatus"
$ws.Range("G2").Style = $ws.Range("F2").Style

# --- Data rows: set Y/N execution status per row ---
$ws.Range("G3").Value = "Y"
$ws.Range("G4").Value = "Y"
$ws.Range("G5").Value = "Y"
$ws.Range("G6").Value = "Y"
$ws.Range("G7").Value = "Y"
$ws.Range("G8").Value = "Y"
$ws.Range("G9").Value = "Y"
$ws.Range("G10").Value = "Y"
$ws.Range("G11").Value = "Y"
$ws.Range("G12").Value = "Y"
$ws.Range("G13").Value = "N"
$ws.Range("G14").Value = "Y"
$ws.Range("G15").Value = "Y"
$ws.Range("G16").Value = "Y"
$ws.Range("G17").Value = "Y"
$ws.Range("G18").Value = "Y"
$ws.Range("G19").Value = "N"
$ws.Range("G20").Value = "N"
$ws.Range("G21").Value = "N"
$ws.Range("G22").Value = "Y"
$ws.Range("G23").Value = "Y"
$ws.Range("G24").Value = "Y"
$ws.Range("G25").Value = "Y"
$ws.Range("G26").Value = "Y"
$ws.Range("G27").Value = "Y"
$ws.Range("G28").Value = "N"

$dataRange = $ws.Range("G3:G28")
$dataRange.Font.Bold = $true
$dataRange.Font.Size = 12
$dataRange.HorizontalAlignment = -4108

# --- Row 28 Keyword column now has "close" ---
$ws.Range("C28").Value = "close"

# --- Column widths ---
$ws.Columns.Item("E").ColumnWidth = 34.59
$ws.Columns.Item("F").ColumnWidth = 7.59
$ws.Columns.Item("G").ColumnWidth = 15.31
$ws.Columns.Item("H").ColumnWidth = 11.45

# --- Row heights (visual adaptation to bold 12pt font in new column) ---
for ($r = 3; $r -le 28; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}
$ws.Rows.Item(14).RowHeight = 135

# --- View / selection ---
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("G20").Select()

$wb.Windows.Item(1).WindowState = -4143
